# Merge the split "<id>" / "p064v_N" / "</id>" runs back into a single
# run per occurrence, e.g. "<id>" + "p064v_1" + "</id>" -> "<id>p064v_1</id>".
# Word's Find/Replace keeps the formatting of the first matched run when it
# collapses a multi-run match into the replacement text, which is exactly
# the Courier New / color 7f6000 formatting the surrounding "<id>" and
# "</id>" runs already have.

$d = $word.ActiveDocument

$ids = @("p064v_1", "p064v_2", "p064v_3", "p064v_4", "p064v_5")

foreach ($id in $ids) {
    $search = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($search, $false, $false, $false, $false, $false, $true, 1, $false, $search, 2)
}
